$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lornoxicam")

# Row 2 (RxNorm / RxCUI) Code column: replace text "C0055477" with numeric value 20890
$ws.Range("D2").Value = 20890

# Update the selected cell shown in the saved file (B19 -> C13)
$ws.Range("C13").Select()
